# Weekly update: insert a new week of Frutilla price data (rows 896-898) at the
# top of the "Terminal La Palmera de La Serena" data block, pushing the rest of
# the historical rows down by three, and re-append the displaced trailing week
# (old rows 951-953) at the new bottom of the range (rows 954-956).
#
# Columns A, B, C, E, F, G, H, I, J, K, Q, R, T are constant for every row in
# this block, so they only need to be (re)written for the three brand-new rows
# that extend the sheet (954-956). Columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado) and
# S (Precio $/Kg) are rewritten for every affected row (896-956).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- constants shared by every row in this block ---------------------------
$colA = 8
$colB = "Terminal La Palmera de La Serena"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100112025
$colJ = "Frutilla"
$colK = "Sin especificar"
$colQ = "`$/bandeja 7 kilos"
$colR = "Provincia de Melipilla"
$colT = 7

# Reference style (date number format) shared by every cell in column D.
$dateFormat = $ws.Range("D2").NumberFormat

# Fill in the constant columns for the three brand-new rows (954-956); rows
# 896-953 already carry these values from before the edit.
for ($r = 954; $r -le 956; $r++) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
    $ws.Cells.Item($r, 11).Value = $colK
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $colR
    $ws.Cells.Item($r, 20).Value = $colT
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
}

# --- per-row data: row, Fecha(serial), Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$rows = @(
    @(896, 44931, "Especial", 500, 12000, 13000, 12500, 1786),
    @(897, 44931, "Primera", 500, 10000, 11000, 10500, 1500),
    @(898, 44931, "Segunda", 500, 8000, 9000, 8500, 1214),
    @(899, 44819, "Especial", 400, 25000, 26000, 25500, 3643),
    @(900, 44819, "Primera", 500, 20000, 21000, 20500, 2929),
    @(901, 44819, "Segunda", 300, 15000, 16000, 15500, 2214),
    @(902, 44295, "Especial", 240, 14500, 15000, 14750, 2107),
    @(903, 44295, "Primera", 300, 12500, 13000, 12750, 1821),
    @(904, 44295, "Segunda", 240, 10500, 11000, 10750, 1536),
    @(905, 44291, "Especial", 240, 14500, 15000, 14750, 2107),
    @(906, 44291, "Primera", 240, 12500, 13000, 12750, 1821),
    @(907, 44291, "Segunda", 200, 10500, 11000, 10750, 1536),
    @(908, 44613, "Especial", 240, 12500, 13000, 12750, 1821),
    @(909, 44613, "Primera", 400, 10500, 11000, 10750, 1536),
    @(910, 44613, "Segunda", 400, 8500, 9000, 8750, 1250),
    @(911, 44630, "Especial", 400, 12500, 13000, 12750, 1821),
    @(912, 44630, "Primera", 400, 10500, 11000, 10750, 1536),
    @(913, 44630, "Segunda", 400, 8500, 9000, 8750, 1250),
    @(914, 44623, "Especial", 400, 12500, 13000, 12750, 1821),
    @(915, 44623, "Primera", 400, 10500, 11000, 10750, 1536),
    @(916, 44623, "Segunda", 300, 8500, 9000, 8750, 1250),
    @(917, 44679, "Especial", 400, 15500, 16000, 15750, 2250),
    @(918, 44679, "Primera", 400, 13500, 14000, 13750, 1964),
    @(919, 44679, "Segunda", 300, 11500, 12000, 11750, 1679),
    @(920, 44414, "Especial", 200, 26500, 27000, 26750, 3821),
    @(921, 44414, "Primera", 240, 21500, 22000, 21750, 3107),
    @(922, 44414, "Segunda", 160, 17500, 18000, 17750, 2536),
    @(923, 44741, "Primera", 320, 22000, 23000, 22500, 3214),
    @(924, 44741, "Segunda", 200, 18000, 19000, 18500, 2643),
    @(925, 44350, "Especial", 400, 18500, 19000, 18750, 2679),
    @(926, 44350, "Primera", 300, 15500, 16000, 15750, 2250),
    @(927, 44350, "Segunda", 240, 12500, 13000, 12750, 1821),
    @(928, 44245, "Especial", 360, 15500, 16000, 15750, 2250),
    @(929, 44245, "Primera", 240, 13500, 14000, 13750, 1964),
    @(930, 44245, "Segunda", 240, 11500, 12000, 11750, 1679),
    @(931, 44565, "Especial", 400, 12000, 13000, 12500, 1786),
    @(932, 44565, "Primera", 300, 10000, 11000, 10500, 1500),
    @(933, 44565, "Segunda", 300, 8000, 9000, 8500, 1214),
    @(934, 44809, "Especial", 300, 25000, 26000, 25500, 3643),
    @(935, 44809, "Primera", 240, 20000, 21000, 20500, 2929),
    @(936, 44809, "Segunda", 300, 15000, 16000, 15500, 2214),
    @(937, 44589, "Especial", 400, 11500, 12000, 11750, 1679),
    @(938, 44589, "Primera", 360, 9500, 10000, 9750, 1393),
    @(939, 44589, "Segunda", 360, 7500, 8000, 7750, 1107),
    @(940, 44622, "Especial", 400, 12500, 13000, 12750, 1821),
    @(941, 44622, "Primera", 400, 10500, 11000, 10750, 1536),
    @(942, 44622, "Segunda", 360, 8500, 9000, 8750, 1250),
    @(943, 44449, "Primera", 300, 26000, 27000, 26500, 3786),
    @(944, 44449, "Segunda", 200, 20000, 21000, 20500, 2929),
    @(945, 44806, "Especial", 240, 24000, 25000, 24500, 3500),
    @(946, 44806, "Primera", 300, 19000, 20000, 19500, 2786),
    @(947, 44806, "Segunda", 300, 15000, 16000, 15500, 2214),
    @(948, 44357, "Especial", 200, 18500, 19000, 18750, 2679),
    @(949, 44357, "Primera", 200, 16500, 17000, 16750, 2393),
    @(950, 44357, "Segunda", 160, 12500, 13000, 12750, 1821),
    @(951, 44911, "Especial", 500, 12000, 13000, 12500, 1786),
    @(952, 44911, "Primera", 600, 10000, 11000, 10500, 1500),
    @(953, 44911, "Segunda", 400, 8000, 9000, 8500, 1214),
    @(954, 44736, "Especial", 300, 24000, 25000, 24500, 3500),
    @(955, 44736, "Primera", 300, 21000, 22000, 21500, 3071),
    @(956, 44736, "Segunda", 300, 16000, 17000, 16500, 2357)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]    # D - Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[7]   # S - Precio $/Kg
}

Write-Output "Applied weekly Frutilla update: rows 896-956 written."
